$d = $word.ActiveDocument

$pairs = @(
    @("983÷5=196, 3", "533÷3=177, 2"),
    @("399÷2=199, 1", "233÷2=116, 1"),
    @("560÷9=62, 2", "647÷9=71, 8"),
    @("917÷8=114, 5", "168÷7=24, 0"),
    @("644÷8=80, 4", "603÷7=86, 1"),
    @("985÷5=197, 0", "888÷9=98, 6"),
    @("500÷4=125, 0", "695÷7=99, 2"),
    @("874÷2=437, 0", "192÷4=48, 0"),
    @("934÷2=467, 0", "686÷7=98, 0"),
    @("322÷5=64, 2", "799÷9=88, 7"),
    @("975÷3=325, 0", "581÷6=96, 5"),
    @("504÷8=63, 0", "780÷5=156, 0"),
    @("319÷5=63, 4", "935÷6=155, 5"),
    @("776÷7=110, 6", "591÷7=84, 3"),
    @("627÷2=313, 1", "873÷3=291, 0"),
    @("179÷6=29, 5", "877÷5=175, 2"),
    @("757÷9=84, 1", "814÷7=116, 2"),
    @("876÷5=175, 1", "923÷3=307, 2"),
    @("930÷3=310, 0", "728÷9=80, 8"),
    @("417÷2=208, 1", "833÷4=208, 1"),
    @("798÷5=159, 3", "937÷4=234, 1"),
    @("860÷3=286, 2", "358÷9=39, 7"),
    @("733÷8=91, 5", "256÷5=51, 1"),
    @("280÷3=93, 1", "164÷7=23, 3"),
    @("965÷9=107, 2", "316÷5=63, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
